$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.454.73"
$ws.Range("E2").Value = "  -0.90%  "

$ws.Range("D3").Value = "2.422.11"
$ws.Range("E3").Value = "  -2.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.02"
$ws.Range("E4").Value = "  +1.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.62"
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.38"
$ws.Range("E6").Value = "  -6.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").Value = "  -4.74%  "

$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("E9").Value = "  -6.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "31.37"
$ws.Range("E10").Value = "  -8.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("E11").Value = "  -3.78%  "

$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").Value = "2.801.77"
$ws.Range("E13").Value = "  -2.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.62"
$ws.Range("E14").Value = "  -5.81%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.77"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.379.73"
$ws.Range("E16").Value = "  -5.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.745"
$ws.Range("E17").Value = "  -5.69%  "

$ws.Range("D18").Value = "41.015.22"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -4.66%  "

$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -3.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.15"
$ws.Range("E21").Value = "  -0.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.52"
$ws.Range("E22").Value = "  -10.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.54"
$ws.Range("E23").Value = "  -3.76%  "

$ws.Range("E24").Value = "  -5.61%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("E26").Value = "  -6.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.17"
$ws.Range("E27").Value = "  -6.24%  "

$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.33"
$ws.Range("E29").Value = "  -4.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.05"
$ws.Range("E30").Value = "  -7.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.03"
$ws.Range("E31").Value = "  -2.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.18"
$ws.Range("E32").Value = "  -8.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.52"
$ws.Range("E33").Value = "  -3.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.49"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0730"
$ws.Range("E35").Value = "  -3.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.01"
$ws.Range("E36").Value = "  -1.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.84"
$ws.Range("E37").Value = "  -5.89%  "

$ws.Range("E38").Value = "  -7.69%  "

$ws.Range("E39").Value = "  -4.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0968"
$ws.Range("E40").Value = "  -8.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.92"
$ws.Range("E41").Value = "  -3.40%  "

$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.66"
$ws.Range("E43").Value = "  -12.22%  "

$ws.Range("D44").Value = "1.902.62"
$ws.Range("E44").Value = "  -4.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0272"
$ws.Range("E45").Value = "  -5.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("E46").Value = "  -9.26%  "

$ws.Range("D47").Value = "2.672.82"
$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("E48").Value = "  -3.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.27"
$ws.Range("E49").Value = "  -5.87%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.97"
$ws.Range("E51").Value = "  -8.77%  "

